# Tender-file tech-part deck: insert a new "file templates" slide after the
# system-architecture slide (position 10), reusing the picture-fill shape
# already present there so the embedded image relationship stays intact.

$p = $ppt.ActivePresentation

# --- 1. Duplicate slide 10 (系统架构图) so the new slide inherits a working
#        image relationship (rId -> ../media/image1.png) ---------------------
$srcSlide = $p.Slides.Item(10)
$range = $srcSlide.Duplicate()
$newSlide = $range.Item(1)

# --- 2. Grab a handle on the template picture shape ("文件系统") before we
#        start deleting things, and make two more copies of it -------------
$template = $null
for ($i = 1; $i -le $newSlide.Shapes.Count; $i++) {
    $shp = $newSlide.Shapes.Item($i)
    if ($shp.Name -eq "文件系统") {
        $template = $shp
    }
}

$fileShape1 = $template
$fileShape2 = $template.Duplicate()
$fileShape3 = $template.Duplicate()

# --- 3. Position + label the three file shapes ------------------------------
$EMU = 12700.0

$fileShape1.Left = 508000 / $EMU
$fileShape1.Top = 3023119 / $EMU
$fileShape1.Width = 2471589 / $EMU
$fileShape1.Height = 689382 / $EMU
$fileShape1.TextFrame.TextRange.Text = "Base.html"

$fileShape2.Left = 4187371 / $EMU
$fileShape2.Top = 3023119 / $EMU
$fileShape2.Width = 2471589 / $EMU
$fileShape2.Height = 689382 / $EMU
$fileShape2.TextFrame.TextRange.Text = "login.html"

$fileShape3.Left = 4187371 / $EMU
$fileShape3.Top = 4081839 / $EMU
$fileShape3.Width = 2471589 / $EMU
$fileShape3.Height = 689382 / $EMU
$fileShape3.TextFrame.TextRange.Text = "logout.html"

# --- 4. Clear the title placeholder + delete every other inherited shape ----
for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $newSlide.Shapes.Item($i)
    if ($shp.Type -eq 14) {
        $shp.TextFrame.TextRange.Text = ""
    } elseif ($shp.Name -ne "文件系统") {
        $shp.Delete()
    }
}

Write-Output $p.Slides.Count
